# fix efforts document and typos in rasd
#
# The single long "Sara" row that recorded three merged dates
# (29/10/2019 + 30/10/2019 + 01/11/2019) with one combined topic string and a
# total of 9 hours is split into three separate dated rows (29/10, 30/10,
# 01/11) each with its own topic and hour count (4 + 3 + 2 = 9, unchanged
# total). The running SUM() formula/range grows to cover the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Row 10 (Matteo table) - "Use case diagrams" topic text is unchanged;
#    nothing to do here content-wise (shared-string index renumbering is
#    an internal bookkeeping detail the engine handles automatically).
# ------------------------------------------------------------------

# ------------------------------------------------------------------
# 2) Rebuild the tail of the "Sara" table (rows 33-34 -> rows 33-46)
# ------------------------------------------------------------------

# Row 33 keeps its existing (blue, fully-bordered) formatting; only the
# content changes: a real date instead of the free-text date label, and a
# shorter topic string.
$ws.Range("A33").Value = 43767
$ws.Range("B33").Value = "Doc structure + Requirements"
$ws.Range("C33").Value = 4
$ws.Rows.Item(33).RowHeight = 29

# Row 34 / Row 35: two brand-new rows, built by duplicating the formatting
# of row 33 (same blue fill / fonts / number format) and then re-shaping
# the borders so the three data rows (33-35) visually form one block.
$ws.Range("A33:C33").Copy()
$ws.Range("A34:C34").PasteSpecial(-4122)
$ws.Range("A33:C33").Copy()
$ws.Range("A35:C35").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A34").Value = 43768
$ws.Range("B34").Value = "Domain assumptions +  Software Sys Attr"
$ws.Range("C34").Value = 3

$ws.Range("A35").Value = 43770
$ws.Range("B35").Value = "State diagrams + User characteristics"
$ws.Range("C35").Value = 2

$ws.Rows.Item(34).RowHeight = 29
$ws.Rows.Item(35).RowHeight = 29

# Remove the bottom border on row 34's A:B cells so rows 33-35 read as one
# continuous bordered block, and drop the left border on C35 to match.
$ws.Range("A34:B34").Borders.Item(9).LineStyle = -4142
$ws.Range("C35").Borders.Item(7).LineStyle = -4142

# ------------------------------------------------------------------
# 3) Ten blank, blue-filled filler rows (36-45), then the relocated total
#    row (46) that now sums the whole expanded range.
# ------------------------------------------------------------------
$ws.Range("A33:C33").Copy()
for ($r = 36; $r -le 45; $r++) {
    $ws.Range("A$r`:C$r").PasteSpecial(-4122)
    $ws.Cells.Item($r, 1).Value = ""
    $ws.Cells.Item($r, 2).Value = ""
    $ws.Cells.Item($r, 3).Value = ""
}
$excel.CutCopyMode = $false

# Old total row (34) moves down to row 46, still green-filled/bold, and its
# SUM now spans the whole widened data range C28:C45.
$ws.Range("A34:C34").Copy()
$ws.Range("A46:C46").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A46").Value = ""
$ws.Range("B46").Value = "Total effort"
$ws.Range("C46").Formula = "=SUM(C28:C45)"

# Clear the old row 34 total's leftover formula/value now that row 34 holds
# new data instead (its format was already refreshed above).
$ws.Range("B34").Formula = ""
$ws.Range("B34").Value = "Domain assumptions +  Software Sys Attr"

# ------------------------------------------------------------------
# 4) Small row-height tweaks elsewhere on the sheet
# ------------------------------------------------------------------
$ws.Rows.Item(4).RowHeight = 29
$ws.Rows.Item(5).RowHeight = 29
$ws.Rows.Item(8).RowHeight = 14.5
$ws.Rows.Item(17).RowHeight = 29
$ws.Rows.Item(13).RowHeight = 27.5
$ws.Rows.Item(25).RowHeight = 27.5

# ------------------------------------------------------------------
# 5) Cosmetic view state (selection / scroll position) to mirror the
#    author's final on-screen state while editing the Sara table.
# ------------------------------------------------------------------
$ws.Range("G35").Select()
$excel.ActiveWindow.ScrollRow = 22

# Make sure every formula (in particular the relocated running total) is
# recomputed against the final values written above.
$excel.Calculate()
